$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new review row (row 25) with the same appid/keyword/email/recovery/time
# pattern as the existing "taxi" rows, a new review comment and a "confirm" status.
$ws.Range("A25").Value = "com.singleton.strechy"
$ws.Range("B25").Value = "taxi"
$ws.Range("C25").Value = "nitanoren23@gmail.com"
$ws.Range("D25").Value = "ronoren61@gmail.com"
$ws.Range("E25").Value = "27/5/2019 15:59"
$ws.Range("F25").Value = "you have a lot of games in the world but this taxi game app is really creative"
$ws.Range("G25").Value = "confirm"

# Add the mailto hyperlinks for the email / recovery-email columns of the new row.
$ws.Hyperlinks.Add($ws.Range("C25"), "mailto:nitanoren23@gmail.com", "", "", "nitanoren23@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D25"), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com")

# Copy the formatting (styles) from row 5, which has the identical column layout,
# onto the new row 25 (applied after the hyperlinks so the normal cell style wins
# over the automatic "Hyperlink" style).
$ws.Range("A5:G5").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)

# Update the active selection to match the edited workbook (cell D25 selected).
$ws.Range("D25").Select()
